$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: status text changed from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both languages / both files.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (F) / "Latest Handback File"
#    (G) columns and the real handback timestamp (H), now that the handback
#    report has been generated. Hyperlinks are rebuilt from scratch so the
#    final relationship / hyperlink order matches row order.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Hyperlinks.Delete()

$wsZhCn.Range("F2").Value = "866128af-c827-49fa-9607-10ad2fd6e7e6.md"
$wsZhCn.Range("G2").Value = "866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-03-12 08:11:21"

$wsZhCn.Range("F3").Value = "ce89a105-6722-4518-b1d2-d1cb6d9401cf.md"
$wsZhCn.Range("G3").Value = "ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-03-12 08:11:21"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/866128af-c827-49fa-9607-10ad2fd6e7e6.md", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/866128af-c827-49fa-9607-10ad2fd6e7e6.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42ee1396c2636bdaf5a46c6c1a7138cced4bafdb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.zh-cn.xlf", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/866128af-c827-49fa-9607-10ad2fd6e7e6.md", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42ee1396c2636bdaf5a46c6c1a7138cced4bafdb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.zh-cn.xlf", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/ce89a105-6722-4518-b1d2-d1cb6d9401cf.md", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/ce89a105-6722-4518-b1d2-d1cb6d9401cf.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42ee1396c2636bdaf5a46c6c1a7138cced4bafdb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.zh-cn.xlf", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/ce89a105-6722-4518-b1d2-d1cb6d9401cf.md", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42ee1396c2636bdaf5a46c6c1a7138cced4bafdb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.zh-cn.xlf", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.zh-cn.xlf")

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of change as zh-cn, different handback time and
#    different xlf hyperlink targets.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Hyperlinks.Delete()

$wsDeDe.Range("F2").Value = "866128af-c827-49fa-9607-10ad2fd6e7e6.md"
$wsDeDe.Range("G2").Value = "866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-03-12 08:11:27"

$wsDeDe.Range("F3").Value = "ce89a105-6722-4518-b1d2-d1cb6d9401cf.md"
$wsDeDe.Range("G3").Value = "ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-03-12 08:11:27"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/866128af-c827-49fa-9607-10ad2fd6e7e6.md", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/866128af-c827-49fa-9607-10ad2fd6e7e6.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc9f80c839625ed6a3f486cd7f0e95d213f53524/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.de-de.xlf", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/866128af-c827-49fa-9607-10ad2fd6e7e6.md", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc9f80c839625ed6a3f486cd7f0e95d213f53524/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.de-de.xlf", "", "", "866128af-c827-49fa-9607-10ad2fd6e7e6.683b67b04efafbb7531982876146a9a139eb7dba.de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/ce89a105-6722-4518-b1d2-d1cb6d9401cf.md", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/ce89a105-6722-4518-b1d2-d1cb6d9401cf.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc9f80c839625ed6a3f486cd7f0e95d213f53524/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.de-de.xlf", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/4fba8ee991953cf936f90a906588a2a9cede023a/e2e/ce89a105-6722-4518-b1d2-d1cb6d9401cf.md", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc9f80c839625ed6a3f486cd7f0e95d213f53524/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.de-de.xlf", "", "", "ce89a105-6722-4518-b1d2-d1cb6d9401cf.08dfde45062107af355bcfe9974de2fe26abb4eb.de-de.xlf")
